$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import")

# Copy the formatting of the last existing data row (286) down across the
# new rows (287-297) so the new cells pick up the same style (s="1") as
# the rest of the translation table.
$ws.Range("A286:C286").Copy()
$ws.Range("A287:C297").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New translation rows to append to the "Import" sheet.
$newRows = @(
    @("cs", "error.Duplicate entry [z_build_name_unique] of [z_build].", "Tento název buildu jste již použili; zvolte prosím jiný."),
    @("cs", "lab.build.table.name", "Jméno buildu"),
    @("cs", "lab.build.table.atomizer", "Atomizér"),
    @("cs", "lab.build.table.cotton", "Vata"),
    @("cs", "lab.build.table.coil", "Spirálka"),
    @("cs", "lab.build.table.ohm", "Odpor buildu"),
    @("cs", "lab.build.table.coils", "Počet spirálek"),
    @("cs", "lab.build.table.created", "Vytvořen"),
    @("cs", "lab.build.created.message", "Build [{{data.name}}] byl uložen."),
    @("cs", "lab.build.table.coilOffset", "Pozice spirálky"),
    @("cs", "lab.build.table.cottonOffset", "Množství vaty")
)

$r = 287
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# Match the saved view state from the diff: scrolled so row 275 is at the
# top, with B297 (the last new row's label cell) selected.
$ws.Activate()
$ws.Range("B297").Select()
$excel.ActiveWindow.ScrollRow = 275
